# Add files via upload
# Rewrites the stock-screener table (columns A-F, rows 2-18) with updated
# ticker lists. Column A keeps its running index (0,1,2,...) with the same
# bordered/bold/centered style already used on A2:A5. Columns C and E end
# up completely empty for every data row, while B, D and F get new ticker
# values (some rows only populate D/F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend column A's index/style down through row 18 -------------------
# Copy the existing formatted cell (A2) and paste its formatting onto the
# new rows so they reuse the same style index instead of creating new ones.
$ws.Range("A2").Copy()
$ws.Range("A6:A18").PasteSpecial(-4122)

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Clear out columns C and E entirely for all data rows ----------------
$ws.Range("C2:C18").ClearContents()
$ws.Range("E2:E18").ClearContents()

# --- Column B (ticker symbols, only some rows populated) -----------------
$ws.Range("B2").Value = "NSE:AVANTIFEED"
$ws.Range("B3").Value = "NSE:BVCL"
$ws.Range("B4").Value = "NSE:DEVYANI"
$ws.Range("B5").Value = "NSE:HINDALCO"
$ws.Range("B6").Value = "NSE:KDDL"
$ws.Range("B7").Value = "NSE:NAUKRI"
$ws.Range("B8:B18").ClearContents()

# --- Column D (ticker symbols, populated through row 18) -----------------
$ws.Range("D2").Value = "NSE:ABCAPITAL"
$ws.Range("D3").Value = "NSE:APLAPOLLO"
$ws.Range("D4").Value = "NSE:BPCL"
$ws.Range("D5").Value = "NSE:DEEPAKNTR"
$ws.Range("D6").Value = "NSE:DIVISLAB"
$ws.Range("D7").Value = "NSE:IEX"
$ws.Range("D8").Value = "NSE:INDHOTEL"
$ws.Range("D9").Value = "NSE:INDIGO"
$ws.Range("D10").Value = "NSE:JSL"
$ws.Range("D11").Value = "NSE:JSWENERGY"
$ws.Range("D12").Value = "NSE:LICI"
$ws.Range("D13").Value = "NSE:MUTHOOTFIN"
$ws.Range("D14").Value = "NSE:NBCC"
$ws.Range("D15").Value = "NSE:OBEROIRLTY"
$ws.Range("D16").Value = "NSE:PAYTM"
$ws.Range("D17").Value = "NSE:POLYCAB"
$ws.Range("D18").Value = "NSE:PRESTIGE"

# --- Column F (ticker symbols, only some rows populated) -----------------
$ws.Range("F2").Value = "NSE:ANGELONE"
$ws.Range("F3").Value = "NSE:BPCL"
$ws.Range("F4").Value = "NSE:INDIGO"
$ws.Range("F5").Value = "NSE:MARICO"
$ws.Range("F6").Value = "NSE:NAUKRI"
$ws.Range("F7").Value = "NSE:NMDC"
$ws.Range("F8").Value = "NSE:PETRONET"
$ws.Range("F9:F18").ClearContents()
